$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 112
$ws.Range("H112").Value = 38463376
$ws.Range("I112").Value = 142857620
$ws.Range("J112").Value = 2341.8948
$ws.Range("K112").Value = 428572860
$ws.Range("L112").Value = 7025.6844
$ws.Range("M112").Value = -428571752
$ws.Range("N112").Value = -9241.6844

# Row 115
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 1750
$ws.Range("J115").Value = 2500
$ws.Range("K115").Value = 5250
$ws.Range("L115").Value = 7500
$ws.Range("M115").Value = -3683

# Row 118
$ws.Range("H118").Value = 1054.5333
$ws.Range("I118").Value = 301.25
$ws.Range("J118").Value = 1915.4286
$ws.Range("K118").Value = 903.75
$ws.Range("L118").Value = 5746.2858
$ws.Range("M118").Value = 753.25
$ws.Range("N118").Value = -9060.2858

# Row 125
$ws.Range("H125").Value = 1840
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1840
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 16560
$ws.Range("N125").Value = -21480
$ws.Range("M125").ClearContents()

# Row 129
$ws.Range("H129").Value = 596.9
$ws.Range("I129").Value = 456.5
$ws.Range("J129").Value = 1158.5
$ws.Range("K129").Value = 1369.5
$ws.Range("L129").Value = 3475.5
$ws.Range("M129").Value = 3630.5
$ws.Range("N129").Value = -13475.5

# Row 132
$ws.Range("H132").Value = 1487898.1
$ws.Range("I132").Value = 2549.5715
$ws.Range("J132").Value = 4087258
$ws.Range("K132").Value = 7648.7145
$ws.Range("L132").Value = 12261774
$ws.Range("M132").Value = -5118.7145
$ws.Range("N132").Value = -12266834

# Row 137
$ws.Range("H137").Value = 3450013.5
$ws.Range("I137").Value = 4001287.8
$ws.Range("J137").Value = 4549.75
$ws.Range("K137").Value = 12003863.4
$ws.Range("L137").Value = 13649.25
$ws.Range("M137").Value = -12001313.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 13158182
$ws.Range("I5").Value = 13158182
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 13158182
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -13158070
$ws.Range("N5").ClearContents()

# Row 45
$ws.Range("H45").Value = 1093.75
$ws.Range("I45").Value = 944.8
$ws.Range("J45").Value = 1342
$ws.Range("K45").Value = 944.8
$ws.Range("L45").Value = 1342
$ws.Range("M45").Value = -567.8
$ws.Range("N45").Value = -2096

# Row 61
$ws.Range("H61").Value = 83502510
$ws.Range("I61").Value = 166834610
$ws.Range("J61").Value = 170416.67
$ws.Range("K61").Value = 166834610
$ws.Range("L61").Value = 170416.67
$ws.Range("M61").Value = -166834398
$ws.Range("N61").Value = -170840.67

# Row 74
$ws.Range("H74").Value = 14581679
$ws.Range("I74").Value = 19667986
$ws.Range("J74").Value = 170475
$ws.Range("K74").Value = 19667986
$ws.Range("L74").Value = 170475
$ws.Range("M74").Value = -19667112
$ws.Range("N74").Value = -172223

# Row 77
$ws.Range("H77").Value = 14581679
$ws.Range("I77").Value = 19667986
$ws.Range("J77").Value = 170475
$ws.Range("K77").Value = 98339930
$ws.Range("L77").Value = 852375
$ws.Range("M77").Value = -98335562
$ws.Range("N77").Value = -861111

# Row 107
$ws.Range("H107").Value = 19296
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 19296
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 19296
$ws.Range("N107").Value = -26976

# Row 109
$ws.Range("H109").Value = 52928.332
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 52928.332
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 52928.332
$ws.Range("N109").Value = -55702.332

# Row 136
$ws.Range("H136").Value = 83502510
$ws.Range("I136").Value = 166834610
$ws.Range("J136").Value = 170416.67
$ws.Range("K136").Value = 500503830
$ws.Range("L136").Value = 511250.01
$ws.Range("M136").Value = -500501280
$ws.Range("N136").Value = -516350.01

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 13158182
$ws.Range("I4").Value = 13158182
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 13158182
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -13158067
$ws.Range("N4").ClearContents()

# Row 103
$ws.Range("H103").Value = 45657
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 45657
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 45657
$ws.Range("N103").Value = -48001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3333.9546
$ws.Range("I31").Value = 1954.5238
$ws.Range("J31").Value = 4593.4346
$ws.Range("K31").Value = 1954.5238
$ws.Range("L31").Value = 4593.4346
$ws.Range("M31").Value = -1659.5238
$ws.Range("N31").Value = -5183.4346

# Row 34
$ws.Range("H34").Value = 3333.9546
$ws.Range("I34").Value = 1954.5238
$ws.Range("J34").Value = 4593.4346
$ws.Range("K34").Value = 1954.5238
$ws.Range("L34").Value = 4593.4346
$ws.Range("M34").Value = -1752.5238
$ws.Range("N34").Value = -4997.4346

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2242.762
$ws.Range("I4").Value = 187
$ws.Range("J4").Value = 2726.4707
$ws.Range("K4").Value = 561
$ws.Range("L4").Value = 8179.4121
$ws.Range("M4").Value = -449
$ws.Range("N4").Value = -8403.4121

# Row 131
$ws.Range("H131").Value = 912.3333
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 934.8
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 2804.4
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -12884.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1185.8823
$ws.Range("I97").Value = 1368.5714
$ws.Range("J97").Value = 333.33334
$ws.Range("K97").Value = 1368.5714
$ws.Range("L97").Value = 333.33334
$ws.Range("M97").Value = -872.5714
$ws.Range("N97").Value = -1325.33334

# Row 102
$ws.Range("H102").Value = 948.7895
$ws.Range("I102").Value = 883.94116
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 883.94116
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 738.05884
$ws.Range("N102").Value = -4744

# Row 122
$ws.Range("H122").Value = 2286.4285
$ws.Range("I122").Value = 1681
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 5043
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -2593
$ws.Range("N122").Value = -16300

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 2833.3333
$ws.Range("I20").Value = 2500
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -2274
$ws.Range("N20").Value = -3452

# Row 136
$ws.Range("H136").Value = 135626.6
$ws.Range("I136").Value = 167967.33
$ws.Range("J136").Value = 114066.11
$ws.Range("K136").Value = 503901.99
$ws.Range("L136").Value = 342198.33
$ws.Range("M136").Value = -501351.99
$ws.Range("N136").Value = -347298.33

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 113422.336
$ws.Range("I132").Value = 113200.445
$ws.Range("J132").Value = 113644.22
$ws.Range("K132").Value = 339601.335
$ws.Range("L132").Value = 340932.66
$ws.Range("M132").Value = -337071.335
$ws.Range("N132").Value = -345992.66
